$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (sorted / corrected) data for rows 2-10, columns A:E
# Modelo, Comparaciones_Significativas, Proporcion_Sig, Mejor_N_Calib, ECRPS_Mejor
$data = @(
    @("AREPD",               "0/10", 0, 60,  0.6709855579037426),
    @("AV-MCPS",              "0/10", 0, 60,  0.6157499770895435),
    @("Block Bootstrapping",  "0/10", 0, 60,  0.6023804666589008),
    @("DeepAR",                "0/10", 0, 20,  0.5562680385510388),
    @("EnCQR-LSTM",            "0/10", 0, 100, 0.8152464016057361),
    @("LSPM",                   "0/10", 0, 100, 0.6292884652770241),
    @("LSPMW",                   "0/10", 0, 60,  0.6483733578521601),
    @("MCPS",                     "0/10", 0, 60,  0.566638548095583),
    @("Sieve Bootstrap",           "0/10", 0, 100, 0.5877694669937914)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
    $ws.Range("C$row").Value = $data[$i][2]
    $ws.Range("D$row").Value = $data[$i][3]
    $ws.Range("E$row").Value = $data[$i][4]
}
